$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear existing sheet contents (old layout used columns A-D ad hoc;
# new layout is a clean two-column A=Japanese / B=English table).
$ws.Cells.Clear()

$ws.Cells.Item(1, 1).Value = '自動実行'
$ws.Cells.Item(1, 2).Value = '自動実行'
$ws.Cells.Item(2, 1).Value = 'リリー'
$ws.Cells.Item(2, 2).Value = 'Lily'
$ws.Cells.Item(3, 1).Value = 'シィナ'
$ws.Cells.Item(3, 2).Value = 'Shina'
$ws.Cells.Item(4, 1).Value = 'アーコイベ用'
$ws.Cells.Item(4, 2).Value = 'アーコイベ用'
$ws.Cells.Item(5, 1).Value = 'ライム'
$ws.Cells.Item(5, 2).Value = 'Lime'
$ws.Cells.Item(6, 1).Value = '自動実行カギ'
$ws.Cells.Item(6, 2).Value = '自動実行カギ'
$ws.Cells.Item(7, 1).Value = '\n<リリー>あ、シィナあんた。
ちゃんと鍵掛けた？'
$ws.Cells.Item(7, 2).Value = '\n<Lily>Oh Shina, you''re back.
Did you lock them up?
'
$ws.Cells.Item(8, 1).Value = '\n<シィナ>は？知らんし。'
$ws.Cells.Item(8, 2).Value = '\n<Shina>Huh, me?
Beats me.'
$ws.Cells.Item(9, 1).Value = '\n<リリー>ちょっと！
しっかりしてよね！
逃げられたらどうすんのよ！'
$ws.Cells.Item(9, 2).Value = '\n<Lily>Whoa whoa whoa!
You need to take this seriously!
What if they got out!'
$ws.Cells.Item(10, 1).Value = '\n<リリー>はぁ～～～・・・
ほんとあんたっていっつも・・・'
$ws.Cells.Item(10, 2).Value = '\n<Lily>Sigh....
Why are you always so... so...'
$ws.Cells.Item(11, 1).Value = '\n<シィナ>はぁ？
知らんにゃ。
鍵持ってないし。アタシじゃねー！'
$ws.Cells.Item(11, 2).Value = '\n<Shina>Whaa?
I dunno-nya.
I don''t have the key. It''s not my fault!'
$ws.Cells.Item(12, 1).Value = '\n<ライム>私も持ってないよ。'
$ws.Cells.Item(12, 2).Value = '\n<Lime>I don''t have a key either.'
$ws.Cells.Item(13, 1).Value = '\n<リリー>ん？'
$ws.Cells.Item(13, 2).Value = '\n<Lily>Hm?'
$ws.Cells.Item(14, 1).Value = '\n<リリー>・・・！！'
$ws.Cells.Item(14, 2).Value = '\n<Lily>...!!'
$ws.Cells.Item(15, 1).Value = '\n<リリー>私が持ってた。
・・・
鍵、閉めたっけ？'
$ws.Cells.Item(15, 2).Value = '\n<Lily>I have the key.
...
Did you at least shut the door?'
$ws.Cells.Item(16, 1).Value = '\n<シィナ>知らんし。
アホにゃ。'
$ws.Cells.Item(16, 2).Value = '\n<Shina>I dunno.
Why don''t you do it yourself-nya.'
$ws.Cells.Item(17, 1).Value = '\n<ライム>とりあえず、確認しに戻ろっか。
リリー。'
$ws.Cells.Item(17, 2).Value = '\n<Lime>Well for now, let''s go back and check,
Lily?'
$ws.Cells.Item(18, 1).Value = '\n<シィナ>ごめんなさいは～？'
$ws.Cells.Item(18, 2).Value = '\n<Shina>Now who should be the one apologizing?~'
$ws.Cells.Item(19, 1).Value = '\n<リリー>えへっ♥'
$ws.Cells.Item(19, 2).Value = '\n<Lily>Eei♥'
$ws.Cells.Item(20, 1).Value = '引っかけ場所'
$ws.Cells.Item(20, 2).Value = '引っかけ場所'
$ws.Cells.Item(21, 1).Value = '<TE:立体起動>'
$ws.Cells.Item(21, 2).Value = '<TE:立体起動>'
$ws.Cells.Item(22, 1).Value = '最初のイベント'
$ws.Cells.Item(22, 2).Value = '最初のイベント'
$ws.Cells.Item(23, 1).Value = '\n<ライム>やっぱりね！'
$ws.Cells.Item(23, 2).Value = '\n<Lime>I knew it!'
$ws.Cells.Item(24, 1).Value = '\n<ライム>逃げる時間そんなにないよなー、って思ってたの。
どこかに隠れてたのかな？
待ってたら来ると思ったよー！'
$ws.Cells.Item(24, 2).Value = '\n<Lime>I was just thinking, "there''s was no time to escape at all-".
You hid somewhere, right?
So I thought to wait here, just in case!'
$ws.Cells.Item(25, 1).Value = 'ーーーーー基本変更点ーーーーー'
$ws.Cells.Item(26, 1).Value = '\n<ライム>にへへへへー♥
つーかまーえた♥'
$ws.Cells.Item(26, 2).Value = '\n<Lime>Ehehehehe-♥
Caught you-♥'
$ws.Cells.Item(27, 1).Value = '\n<\n[3]>ぬるぬるして動きにくいでしょー。
早く逃げないと白いの出させちゃうぞー？'
$ws.Cells.Item(27, 2).Value = '\n<\n[3]>It''s hard to move when everything is so slippery right?
If you don''t escape soon you''ll let out your white stuff you know?'
$ws.Cells.Item(28, 1).Value = '\C[3]※捕まるとタイミングバーが表示されます。
\C[0]タイミングよく黄か赤で止めてください。
赤で止めると被ダメージが半減します。'
$ws.Cells.Item(28, 2).Value = '\C[3]※A timing bar will display if you''re caught.
\C[0]Go ahead and stop in the red or yellow areas.
Red areas will reduce damage taken by half.'
$ws.Cells.Item(29, 1).Value = '\n<\n[3]>むにゅー・・・♥
柔らかくて溶けちゃいそうでしょー♥
気持ちいい気持ちいいー♥'
$ws.Cells.Item(29, 2).Value = '\n<\n[3]>Squish-...♥
They''re so soft you''re about to melt aren''t you?♥
Feels so good- Feels so good-♥'
$ws.Cells.Item(30, 1).Value = '\n<\n[3]>あれー？もう出ちゃうのー？
おっぱい我慢できなかったー？あはは♥
じゃあ一回だけ、出しちゃおっかー♪'
$ws.Cells.Item(30, 2).Value = '\n<\n[3]>What''s this? You''re about to cum already-?
So you couldn''t resist my boobs after all? Ahaha♥
Well then, I''ll give you one more pump, let it all out♪'
$ws.Cells.Item(31, 1).Value = '\n<\n[3]>あっあっ♥おっぱいの間でぴくぴくしてるー♥
おちんちん喜んでくれたみたい♥
うれしー♥'
$ws.Cells.Item(31, 2).Value = '\n<\n[3]>Ahh ahh♥ It''s shooting between my boobies-♥
Your penis seems delighted too♥ I''m so happy-♥'
$ws.Cells.Item(32, 1).Value = '\C[1]SAN値が1下がった・・・（現在SAN値\v[270]）'
$ws.Cells.Item(32, 2).Value = '\C[1]Sanity decreased by 1... (Current Sanity: \v[270]）'
$ws.Cells.Item(33, 1).Value = '\n<ライム>もうー。ひょっとして全然抵抗する気ないのー？
ゲームはまだ始まったばかりなのに。
次は本気で搾っちゃうからねー？'
$ws.Cells.Item(33, 2).Value = '\n<Lime>Whoa- Were you seriously not trying to resist at all-?
And the game has only just started, next time I''ll
milk you for real okay?'
$ws.Cells.Item(34, 1).Value = '\n<ライム>あっ！逃げられちゃった・・・
まぁいいっかー。'
$ws.Cells.Item(34, 2).Value = '\n<Lime>Ah! You escaped...
Oh well-.'
$ws.Cells.Item(35, 1).Value = '\n<ライム>えーっと、リリーの言ってたこと、
聞こえてたかな？
脱出ゲームがどーのこーの・・・'
$ws.Cells.Item(35, 2).Value = '\n<Lime>That is...you heard what Lily said right?
About this escape game thingy...'
$ws.Cells.Item(36, 1).Value = '\n<ライム>この館は広いからねー。
なかなか出られないと思うけど・・・
頑張って色んなところ、探してみてね。'
$ws.Cells.Item(36, 2).Value = '\n<Lime>This mansion is very, very big-.
I don''t think you''ll get out...but feel free to try your hardest,
and search everywhere you can okay?'
$ws.Cells.Item(37, 1).Value = '\n<ライム>リリーもシィナも楽しそうだから。
簡単に諦めたりしちゃーダメだよー？
いっぱい遊ばれてね。'
$ws.Cells.Item(37, 2).Value = '\n<Lime>Lily and Shina seem to be happy after all.
So you shouldn''t give up so easily OK? Enjoy yourself
as much as possible now-'
$ws.Cells.Item(38, 1).Value = '\n<ライム>次は一回射精したぐらいじゃ
許してあげないからねー？くすくす♥
ばいばーい♥'
$ws.Cells.Item(38, 2).Value = '\n<Lime>Next time you cum, I won''t go easy on you now alright?
Teehee♥ Bye bye♥'
$ws.Cells.Item(39, 1).Value = '\n<ライム>じゃー頑張ってねー。
次は私も本気でぴゅっぴゅさせにいっちゃうからねー。'
$ws.Cells.Item(39, 2).Value = '\n<Lime>Now then, do your best okay?
Next time I''ll make you go pew pew for realsies-'
$ws.Cells.Item(40, 1).Value = 'ドア'
$ws.Cells.Item(40, 2).Value = 'ドア'
$ws.Cells.Item(41, 1).Value = '玄関扉は固く閉ざされている・・・'
$ws.Cells.Item(41, 2).Value = '玄関扉は固く閉ざされている・・・'
$ws.Cells.Item(42, 1).Value = '\n<？？？>開かないよ。'
$ws.Cells.Item(42, 2).Value = '\n<？？？>開かないよ。'
$ws.Cells.Item(43, 1).Value = '\n<？？？>あたしも閉じ込められて困ってんの。
あんた誰？'
$ws.Cells.Item(43, 2).Value = '\n<？？？>あたしも閉じ込められて困ってんの。
あんた誰？'
$ws.Cells.Item(44, 1).Value = 'MP_SET_MOVIE e8_Ta'
$ws.Cells.Item(44, 2).Value = 'MP_SET_MOVIE e8_Ta'
$ws.Cells.Item(45, 1).Value = 'MP_SET_LOOP 6 on'
$ws.Cells.Item(45, 2).Value = 'MP_SET_LOOP 6 on'
$ws.Cells.Item(46, 1).Value = '\n<アーコ>あーあー警戒しないで。
あたし中立だから。ちゅーりつ。
知ってる？ちゅーりつって言葉。チューリップじゃないよ。'
$ws.Cells.Item(46, 2).Value = '\n<アーコ>あーあー警戒しないで。
あたし中立だから。ちゅーりつ。
知ってる？ちゅーりつって言葉。チューリップじゃないよ。'
$ws.Cells.Item(47, 1).Value = '\n<？？？>こんなところに人間が居るなんて珍しーね。
リリーたちに捕まったのかな？'
$ws.Cells.Item(47, 2).Value = '\n<？？？>こんなところに人間が居るなんて珍しーね。
リリーたちに捕まったのかな？'
$ws.Cells.Item(48, 1).Value = '\n<？？？>あたしアーコ。
黒猫のアーコだよ。よろしく。
あんたは？'
$ws.Cells.Item(48, 2).Value = '\n<？？？>あたしアーコ。
黒猫のアーコだよ。よろしく。
あんたは？'
$ws.Cells.Item(49, 1).Value = '\n[1]・・・'
$ws.Cells.Item(49, 2).Value = '\n[1]・・・'
$ws.Cells.Item(50, 1).Value = '擬態型かァ～？'
$ws.Cells.Item(50, 2).Value = '擬態型かァ～？'
$ws.Cells.Item(51, 1).Value = '\n<アーコ>友好型だよ。アーコだよ。
名前。教えてよ。'
$ws.Cells.Item(51, 2).Value = '\n<アーコ>友好型だよ。アーコだよ。
名前。教えてよ。'
$ws.Cells.Item(52, 1).Value = '\n<アーコ>\n[1]。
うん。覚えた。フツーだね。
フツーにいい名前。フツー。'
$ws.Cells.Item(52, 2).Value = '\n<アーコ>\n[1]。
うん。覚えた。フツーだね。
フツーにいい名前。フツー。'
$ws.Cells.Item(53, 1).Value = '\n<シィナ>アーコ来てんのかー！？'
$ws.Cells.Item(53, 2).Value = '\n<シィナ>アーコ来てんのかー！？'
$ws.Cells.Item(54, 1).Value = '\n<アーコ>隠れて！'
$ws.Cells.Item(54, 2).Value = '\n<アーコ>隠れて！'
$ws.Cells.Item(55, 1).Value = '\n<シィナ>タバコにゃ。
持ってきたにゃ？'
$ws.Cells.Item(55, 2).Value = '\n<シィナ>タバコにゃ。
持ってきたにゃ？'
$ws.Cells.Item(56, 1).Value = '\n<アーコ>ほいほい。
税金上がったから1箱540Gだよ。
何個欲しいの？しぃにゃん。'
$ws.Cells.Item(56, 2).Value = '\n<アーコ>ほいほい。
税金上がったから1箱540Gだよ。
何個欲しいの？しぃにゃん。'
$ws.Cells.Item(57, 1).Value = '\n<シィナ>淫魔が税金納めてんじゃねー！
カートンにゃ！
あと次しぃにゃんって呼んだら殺す！'
$ws.Cells.Item(57, 2).Value = '\n<シィナ>淫魔が税金納めてんじゃねー！
カートンにゃ！
あと次しぃにゃんって呼んだら殺す！'
$ws.Cells.Item(58, 1).Value = '\n<アーコ>5400G！まいどー♪
あたし、しぃにゃんの怒った顔好き。'
$ws.Cells.Item(58, 2).Value = '\n<アーコ>5400G！まいどー♪
あたし、しぃにゃんの怒った顔好き。'
$ws.Cells.Item(59, 1).Value = '\n<シィナ>アーココロス！'
$ws.Cells.Item(59, 2).Value = '\n<シィナ>アーココロス！'
$ws.Cells.Item(60, 1).Value = '\n<アーコ>（やっぱりシィナは怒った時の顔が一番可愛い・・・）'
$ws.Cells.Item(60, 2).Value = '\n<アーコ>（やっぱりシィナは怒った時の顔が一番可愛い・・・）'
$ws.Cells.Item(61, 1).Value = '\n<シィナ>あーお前と話してるとほんと調子狂うにゃ。
はよ帰れし。
じゃーなバーカ。'
$ws.Cells.Item(61, 2).Value = '\n<シィナ>あーお前と話してるとほんと調子狂うにゃ。
はよ帰れし。
じゃーなバーカ。'
$ws.Cells.Item(62, 1).Value = '\n<アーコ>帰りたくても帰れないよ。
玄関開かないもん。'
$ws.Cells.Item(62, 2).Value = '\n<アーコ>帰りたくても帰れないよ。
玄関開かないもん。'
$ws.Cells.Item(63, 1).Value = '\n<シィナ>あ。そうにゃ。
今人間と脱出ゲームやってるから全部閉め切ってるにゃん。'
$ws.Cells.Item(63, 2).Value = '\n<シィナ>あ。そうにゃ。
今人間と脱出ゲームやってるから全部閉め切ってるにゃん。'
$ws.Cells.Item(64, 1).Value = '\n<アーコ>脱出ゲーム？'
$ws.Cells.Item(64, 2).Value = '\n<アーコ>脱出ゲーム？'
$ws.Cells.Item(65, 1).Value = '\n<シィナ>人間がここから脱出するのが先か
アタシらにチンポ犯されまくって干物になるのが先かの・・・'
$ws.Cells.Item(65, 2).Value = '\n<シィナ>人間がここから脱出するのが先か
アタシらにチンポ犯されまくって干物になるのが先かの・・・'
$ws.Cells.Item(66, 1).Value = '\n<シィナ>まさにセイシをかけたゲームにゃ！'
$ws.Cells.Item(66, 2).Value = '\n<シィナ>まさにセイシをかけたゲームにゃ！'
$ws.Cells.Item(67, 1).Value = '\n<アーコ>そう・・・
（また妙なことやってる・・・）'
$ws.Cells.Item(67, 2).Value = '\n<アーコ>そう・・・
（また妙なことやってる・・・）'
$ws.Cells.Item(68, 1).Value = '\n<シィナ>邪魔すんなよ。'
$ws.Cells.Item(68, 2).Value = '\n<シィナ>邪魔すんなよ。'
$ws.Cells.Item(69, 1).Value = '\n<アーコ>えっ？ちょっと待って。
それってあたしも出られないってコト？'
$ws.Cells.Item(69, 2).Value = '\n<アーコ>えっ？ちょっと待って。
それってあたしも出られないってコト？'
$ws.Cells.Item(70, 1).Value = '\n<シィナ>お前のことは知らんし。
大人しく床の溝でもなぞっとけにゃ。
カートンありがとにゃー。'
$ws.Cells.Item(70, 2).Value = '\n<シィナ>お前のことは知らんし。
大人しく床の溝でもなぞっとけにゃ。
カートンありがとにゃー。'
$ws.Cells.Item(71, 1).Value = '\n<アーコ>なるほどねー。
事情は分かったよ。'
$ws.Cells.Item(71, 2).Value = '\n<アーコ>なるほどねー。
事情は分かったよ。'
$ws.Cells.Item(72, 1).Value = '\n<アーコ>どうやら\n[1]が出口を見つけるか干物になるまで
あたしも出られなくなっちゃったみたい。
うーん。'
$ws.Cells.Item(72, 2).Value = '\n<アーコ>どうやら\n[1]が出口を見つけるか干物になるまで
あたしも出られなくなっちゃったみたい。
うーん。'
$ws.Cells.Item(73, 1).Value = '\n<アーコ>今ここで干物にしちゃうか・\.・\.・\.'
$ws.Cells.Item(73, 2).Value = '\n<アーコ>今ここで干物にしちゃうか・\.・\.・\.'
$ws.Cells.Item(74, 1).Value = '\n<アーコ>冗談だよ。
うける。'
$ws.Cells.Item(74, 2).Value = '\n<アーコ>冗談だよ。
うける。'
$ws.Cells.Item(75, 1).Value = '\n<アーコ>ところであたしは淫魔だけどお店をやってるんだよ。
黒猫屋！
精力剤とか必要でしょ？売ってあげるよ！'
$ws.Cells.Item(75, 2).Value = '\n<アーコ>ところであたしは淫魔だけどお店をやってるんだよ。
黒猫屋！
精力剤とか必要でしょ？売ってあげるよ！'
$ws.Cells.Item(76, 1).Value = '\n<アーコ>お金無いの？
じゃあ、もしこの館でお金になりそうなものあったら
あたしに持ってきて。買い取ってあげる。'
$ws.Cells.Item(76, 2).Value = '\n<アーコ>お金無いの？
じゃあ、もしこの館でお金になりそうなものあったら
あたしに持ってきて。買い取ってあげる。'
$ws.Cells.Item(77, 1).Value = '\n<アーコ>あたしここに居るから。
何かあったら話しかけてね。
じゃ。'
$ws.Cells.Item(77, 2).Value = '\n<アーコ>あたしここに居るから。
何かあったら話しかけてね。
じゃ。'
$ws.Cells.Item(78, 1).Value = 'シィナイベ用'
$ws.Cells.Item(78, 2).Value = 'シィナイベ用'
$ws.Cells.Item(79, 1).Value = 'アーコ'
$ws.Cells.Item(79, 2).Value = 'アーコ'
$ws.Cells.Item(80, 1).Value = '\n<アーコ>やぁやぁ。アーコだよ。
何か欲しい物とかあるの？
見てく？'
$ws.Cells.Item(80, 2).Value = '\n<アーコ>やぁやぁ。アーコだよ。
何か欲しい物とかあるの？
見てく？'
$ws.Cells.Item(81, 1).Value = '\n<アーコ>武器になりそうなものとか見つけた？
まずはそういうの探すといいかもね。'
$ws.Cells.Item(81, 2).Value = '\n<アーコ>武器になりそうなものとか見つけた？
まずはそういうの探すといいかもね。'
$ws.Cells.Item(82, 1).Value = '\n<アーコ>何食べたらそんなにおっぱい大きくなるの？
横縞模様だからおっぱい大きく見えるの？
錯覚？'
$ws.Cells.Item(82, 2).Value = '\n<アーコ>何食べたらそんなにおっぱい大きくなるの？
横縞模様だからおっぱい大きく見えるの？
錯覚？'
$ws.Cells.Item(83, 1).Value = '\n<アーコ>囚人服に合ってるよ。
まさに囚人って感じ。
褒めてないよ。'
$ws.Cells.Item(83, 2).Value = '\n<アーコ>囚人服に合ってるよ。
まさに囚人って感じ。
褒めてないよ。'
$ws.Cells.Item(84, 1).Value = '\n<アーコ>あたしは淫魔だし、気持ちが分かるから
どっちの味方もしないよ。
でも仲良くしたいって思ってる。'
$ws.Cells.Item(84, 2).Value = '\n<アーコ>あたしは淫魔だし、気持ちが分かるから
どっちの味方もしないよ。
でも仲良くしたいって思ってる。'
$ws.Cells.Item(85, 1).Value = '\n<アーコ>シィナは猫の時からの友達だよ。
リリーに淫魔にされたんだってさ。
あたしは違うやつに淫魔にされたけど。'
$ws.Cells.Item(85, 2).Value = '\n<アーコ>シィナは猫の時からの友達だよ。
リリーに淫魔にされたんだってさ。
あたしは違うやつに淫魔にされたけど。'
$ws.Cells.Item(86, 1).Value = '\n<アーコ>あたしも淫魔だから精液が一番美味しいと思うけど
カリカリもかつおぶしも同じくらい美味しいと思う。
一番がいっぱいあるといいね。'
$ws.Cells.Item(86, 2).Value = '\n<アーコ>あたしも淫魔だから精液が一番美味しいと思うけど
カリカリもかつおぶしも同じくらい美味しいと思う。
一番がいっぱいあるといいね。'
$ws.Cells.Item(87, 1).Value = '\n<アーコ>ここは元々大金持ちが住んでたんだって。
メイドも使用人もいっぱい居たってさ。
なんでこんな森の中に？って感じ。'
$ws.Cells.Item(87, 2).Value = '\n<アーコ>ここは元々大金持ちが住んでたんだって。
メイドも使用人もいっぱい居たってさ。
なんでこんな森の中に？って感じ。'
$ws.Cells.Item(88, 1).Value = '\n<アーコ>みんなお金に興味ないから
あたしがこの館から色々持って行っても何とも思わないんだよ。
だから、ここで仕入れて町で売ったりする。'
$ws.Cells.Item(88, 2).Value = '\n<アーコ>みんなお金に興味ないから
あたしがこの館から色々持って行っても何とも思わないんだよ。
だから、ここで仕入れて町で売ったりする。'
$ws.Cells.Item(89, 1).Value = '\n<アーコ>その見返りにタバコとか持ってきてあげるの。
まぁ、それもお金取るけどね。'
$ws.Cells.Item(89, 2).Value = '\n<アーコ>その見返りにタバコとか持ってきてあげるの。
まぁ、それもお金取るけどね。'
$ws.Cells.Item(90, 1).Value = '\n<アーコ>淫魔の巣に閉じ込められるのは
人間にとって辛いでしょ。
少しでもえっちな気分になるとすぐ硬くなっちゃう。'
$ws.Cells.Item(90, 2).Value = '\n<アーコ>淫魔の巣に閉じ込められるのは
人間にとって辛いでしょ。
少しでもえっちな気分になるとすぐ硬くなっちゃう。'
$ws.Cells.Item(91, 1).Value = '\n<アーコ>行き詰ったらとりあえず攻撃してみたらいいよ。
壊せるものとかあるかもしれないし。'
$ws.Cells.Item(91, 2).Value = '\n<アーコ>行き詰ったらとりあえず攻撃してみたらいいよ。
壊せるものとかあるかもしれないし。'
$ws.Cells.Item(92, 1).Value = '\n<アーコ>お兄ちゃんも下に居るの？
面白そうだから後で行ってみよー。'
$ws.Cells.Item(92, 2).Value = '\n<アーコ>お兄ちゃんも下に居るの？
面白そうだから後で行ってみよー。'
$ws.Cells.Item(93, 1).Value = '\n<アーコ>妹も捕まってるの？
後で行ってみよー。
仲良くなれるかな？'
$ws.Cells.Item(93, 2).Value = '\n<アーコ>妹も捕まってるの？
後で行ってみよー。
仲良くなれるかな？'
$ws.Cells.Item(94, 1).Value = '買い物'
$ws.Cells.Item(94, 2).Value = '買い物'
$ws.Cells.Item(95, 1).Value = '用事はない'
$ws.Cells.Item(95, 2).Value = '用事はない'
$ws.Cells.Item(96, 1).Value = '\n<アーコ>何を買ってくれるのかな？'
$ws.Cells.Item(96, 2).Value = '\n<アーコ>何を買ってくれるのかな？'
$ws.Cells.Item(97, 1).Value = '\n<アーコ>ははーん。さては冷やかしだな？
ふーん！'
$ws.Cells.Item(97, 2).Value = '\n<アーコ>ははーん。さては冷やかしだな？
ふーん！'
$ws.Cells.Item(98, 1).Value = '\n<アーコ>まいどあり～！'
$ws.Cells.Item(98, 2).Value = '\n<アーコ>まいどあり～！'
$ws.Cells.Item(99, 1).Value = '\n<アーコ>ふーん。'
$ws.Cells.Item(99, 2).Value = '\n<アーコ>ふーん。'
$ws.Cells.Item(100, 1).Value = '食糧庫ドア'
$ws.Cells.Item(100, 2).Value = '食糧庫ドア'
$ws.Cells.Item(101, 1).Value = '食糧庫の鍵を開けた！'
$ws.Cells.Item(101, 2).Value = '食糧庫の鍵を開けた！'
$ws.Cells.Item(102, 1).Value = '食糧庫のようだ。
鍵がかかっている・・・'
$ws.Cells.Item(102, 2).Value = '食糧庫のようだ。
鍵がかかっている・・・'
$ws.Cells.Item(103, 1).Value = 'いわ'
$ws.Cells.Item(103, 2).Value = 'いわ'
$ws.Cells.Item(104, 1).Value = ' <enemy:145>'
$ws.Cells.Item(104, 2).Value = ' <enemy:145>'
$ws.Cells.Item(105, 1).Value = '瓦礫を壊せば通れそうだ・・・'
$ws.Cells.Item(105, 2).Value = '瓦礫を壊せば通れそうだ・・・'
$ws.Cells.Item(106, 1).Value = 'EV023'
$ws.Cells.Item(106, 2).Value = 'EV023'
$ws.Cells.Item(107, 1).Value = '浴場の鍵を開けた！'
$ws.Cells.Item(107, 2).Value = '浴場の鍵を開けた！'
$ws.Cells.Item(108, 1).Value = 'この先は浴場のようだ。
鍵がかかっている・・・'
$ws.Cells.Item(108, 2).Value = 'この先は浴場のようだ。
鍵がかかっている・・・'
$ws.Cells.Item(109, 1).Value = '応接室の鍵を開けた！'
$ws.Cells.Item(109, 2).Value = '応接室の鍵を開けた！'
$ws.Cells.Item(110, 1).Value = '応接室の扉だ。
鍵がかかっている・・・'
$ws.Cells.Item(110, 2).Value = '応接室の扉だ。
鍵がかかっている・・・'
$ws.Cells.Item(111, 1).Value = '倉庫の鍵を開けた！'
$ws.Cells.Item(111, 2).Value = '倉庫の鍵を開けた！'
$ws.Cells.Item(112, 1).Value = '倉庫の扉だ。
鍵がかかっている・・・'
$ws.Cells.Item(112, 2).Value = '倉庫の扉だ。
鍵がかかっている・・・'
$ws.Cells.Item(113, 1).Value = '梯子'
$ws.Cells.Item(113, 2).Value = '梯子'
$ws.Cells.Item(114, 1).Value = '報告書'
$ws.Cells.Item(114, 2).Value = '報告書'
$ws.Cells.Item(115, 1).Value = '時計が時を刻んでいる・・・'
$ws.Cells.Item(115, 2).Value = '時計が時を刻んでいる・・・'
$ws.Cells.Item(116, 1).Value = '\n<\n[1]>（どれだけの時間、ここにいるんだろう・・・）'
$ws.Cells.Item(116, 2).Value = '\n<\n[1]>（どれだけの時間、ここにいるんだろう・・・）'
$ws.Cells.Item(117, 1).Value = '\n<\n[1]>（時間の感覚がおかしくなる・・・
早くここから抜け出さないと・・・）'
$ws.Cells.Item(117, 2).Value = '\n<\n[1]>（時間の感覚がおかしくなる・・・
早くここから抜け出さないと・・・）'
